$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.172937393188477
$ws.Range("B1").Value = 2.50653076171875
$ws.Range("C1").Value = 2.63999080657959
$ws.Range("D1").Value = 3.250929117202759
$ws.Range("E1").Value = 2.325813770294189
